$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the three new columns (F: Shunt Admittans (ohm), G: Shunt Admittans (p.u.),
# H: Shunt Admittans (p.u.) half), one triple per existing data row.
$shunt = @(
    @("62589.9504j",  "695.4439j",  "347.7219j"),
    @("51383.8894j",  "570.9321j",  "285.4661j"),
    @("130646.2721j", "1451.6252j", "725.8126j"),
    @("42954.9964j",  "477.2777j",  "238.6389j"),
    @("254996.7925j", "2833.2977j", "1416.6488j"),
    @("194829.01j",   "2164.7668j", "1082.3834j"),
    @("135019.3691j", "1500.2152j", "750.1076j"),
    @("289171.0374j", "3213.0115j", "1606.5058j"),
    @("318962.7605j", "3544.0307j", "1772.0153j")
)

# Row 1 - new headers (written first so the new shared strings are appended
# immediately after the existing header strings, matching original workbook order).
$ws.Range("F1").Value2 = "Shunt Admittans (ohm)"
$ws.Range("G1").Value2 = "Shunt Admittans (p.u.)"
$ws.Range("H1").Value2 = "Shunt Admittans (p.u.) half"

# Rows 2-10 - new data, written row by row so each row's new strings are appended
# to the shared string table right after that row is processed.
for ($i = 0; $i -lt $shunt.Length; $i++) {
    $row = $i + 2
    $vals = $shunt[$i]
    $ws.Range("F$row").Value2 = $vals[0]
    $ws.Range("G$row").Value2 = $vals[1]
    $ws.Range("H$row").Value2 = $vals[2]
}

# Column widths for the changed/added columns.
$ws.Columns.Item(2).ColumnWidth = 22.998697916666668
$ws.Columns.Item(3).ColumnWidth = 26.166666666666668
$ws.Columns.Item(4).ColumnWidth = 23.498697916666668
$ws.Columns.Item(5).ColumnWidth = 13.498697916666666
$ws.Columns.Item(6).ColumnWidth = 20.666666666666668
$ws.Columns.Item(7).ColumnWidth = 21.666666666666668
$ws.Columns.Item(8).ColumnWidth = 24.830729166666668

# Update the active selection to match the saved workbook state.
$null = $ws.Range("C9").Select()
